$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-49 (columns B and C) with new computed values
$ws.Cells.Item(2, 2).Value = 1.291824554907194
$ws.Cells.Item(2, 3).Value = 3.798797808720669
$ws.Cells.Item(3, 2).Value = 4.205592449451212
$ws.Cells.Item(3, 3).Value = 7.799089205519206
$ws.Cells.Item(4, 2).Value = 5.751052597704876
$ws.Cells.Item(4, 3).Value = 11.80723893682239
$ws.Cells.Item(5, 2).Value = 8.638708594068996
$ws.Cells.Item(5, 3).Value = 16.36966594639864
$ws.Cells.Item(6, 2).Value = 9.695222951751587
$ws.Cells.Item(6, 3).Value = 20.85660928739086
$ws.Cells.Item(7, 2).Value = 11.61535050135384
$ws.Cells.Item(7, 3).Value = 25.16696213914714
$ws.Cells.Item(8, 2).Value = 12.4584762150205
$ws.Cells.Item(8, 3).Value = 29.47105023982667
$ws.Cells.Item(9, 2).Value = 17.23736076129476
$ws.Cells.Item(9, 3).Value = 33.56079632482943
$ws.Cells.Item(10, 2).Value = 19.19797079814772
$ws.Cells.Item(10, 3).Value = 37.60826945831809
$ws.Cells.Item(11, 2).Value = 20.48604650860038
$ws.Cells.Item(11, 3).Value = 41.57864695295302
$ws.Cells.Item(12, 2).Value = 20.80151027509101
$ws.Cells.Item(12, 3).Value = 45.93346337674929
$ws.Cells.Item(13, 2).Value = 24.41973628527606
$ws.Cells.Item(13, 3).Value = 49.94132166399145
$ws.Cells.Item(14, 2).Value = 26.87939180143032
$ws.Cells.Item(14, 3).Value = 54.02507766530126
$ws.Cells.Item(15, 2).Value = 27.53614525264342
$ws.Cells.Item(15, 3).Value = 58.07849739296265
$ws.Cells.Item(16, 2).Value = 31.58296887728824
$ws.Cells.Item(16, 3).Value = 62.05640333452591
$ws.Cells.Item(17, 2).Value = 34.68704387387611
$ws.Cells.Item(17, 3).Value = 65.91352199827283
$ws.Cells.Item(18, 2).Value = 38.64743820430902
$ws.Cells.Item(18, 3).Value = 69.7720871589727
$ws.Cells.Item(19, 2).Value = 40.55907342470277
$ws.Cells.Item(19, 3).Value = 73.9389129247363
$ws.Cells.Item(20, 2).Value = 42.37943011009924
$ws.Cells.Item(20, 3).Value = 77.70306917175786
$ws.Cells.Item(21, 2).Value = 44.00395111117064
$ws.Cells.Item(21, 3).Value = 81.71319173552008
$ws.Cells.Item(22, 2).Value = 44.59922084790226
$ws.Cells.Item(22, 3).Value = 85.82069074118658
$ws.Cells.Item(23, 2).Value = 46.36562880888329
$ws.Cells.Item(23, 3).Value = 90.20381831842074
$ws.Cells.Item(24, 2).Value = 48.84714893517194
$ws.Cells.Item(24, 3).Value = 94.4012615668664
$ws.Cells.Item(25, 2).Value = 52.34175178936282
$ws.Cells.Item(25, 3).Value = 98.48076229785276
$ws.Cells.Item(26, 2).Value = 54.20255320621909
$ws.Cells.Item(26, 3).Value = 102.9427312318516
$ws.Cells.Item(27, 2).Value = 56.22843476167699
$ws.Cells.Item(27, 3).Value = 106.9981732059428
$ws.Cells.Item(28, 2).Value = 57.47847296070766
$ws.Cells.Item(28, 3).Value = 110.9373131761006
$ws.Cells.Item(29, 2).Value = 58.75910863281327
$ws.Cells.Item(29, 3).Value = 114.7686029287055
$ws.Cells.Item(30, 2).Value = 63.25519729710583
$ws.Cells.Item(30, 3).Value = 118.5241117311727
$ws.Cells.Item(31, 2).Value = 63.70922619699857
$ws.Cells.Item(31, 3).Value = 122.7402153682025
$ws.Cells.Item(32, 2).Value = 64.75464939062627
$ws.Cells.Item(32, 3).Value = 126.6020810154741
$ws.Cells.Item(33, 2).Value = 66.01779775656644
$ws.Cells.Item(33, 3).Value = 130.6071938997788
$ws.Cells.Item(34, 2).Value = 68.97855807201829
$ws.Cells.Item(34, 3).Value = 134.7222737020898
$ws.Cells.Item(35, 2).Value = 70.17208105251981
$ws.Cells.Item(35, 3).Value = 138.4882093071851
$ws.Cells.Item(36, 2).Value = 70.33896580910114
$ws.Cells.Item(36, 3).Value = 142.2409714821126
$ws.Cells.Item(37, 2).Value = 71.56168300018982
$ws.Cells.Item(37, 3).Value = 146.8978061832853
$ws.Cells.Item(38, 2).Value = 72.81298911379822
$ws.Cells.Item(38, 3).Value = 151.0626736888793
$ws.Cells.Item(39, 2).Value = 74.79068106688069
$ws.Cells.Item(39, 3).Value = 154.8507455488892
$ws.Cells.Item(40, 2).Value = 76.45065630195799
$ws.Cells.Item(40, 3).Value = 159.6130097113086
$ws.Cells.Item(41, 2).Value = 78.82412647805121
$ws.Cells.Item(41, 3).Value = 163.3652432026426
$ws.Cells.Item(42, 2).Value = 80.14467163744452
$ws.Cells.Item(42, 3).Value = 167.8412029911188
$ws.Cells.Item(43, 2).Value = 80.78554062555102
$ws.Cells.Item(43, 3).Value = 171.8263159837037
$ws.Cells.Item(44, 2).Value = 82.56064849314504
$ws.Cells.Item(44, 3).Value = 175.8922167608311
$ws.Cells.Item(45, 2).Value = 86.35548680370798
$ws.Cells.Item(45, 3).Value = 179.8972974896811
$ws.Cells.Item(46, 2).Value = 88.8724730799278
$ws.Cells.Item(46, 3).Value = 183.776697798123
$ws.Cells.Item(47, 2).Value = 91.45246408717054
$ws.Cells.Item(47, 3).Value = 188.323730458367
$ws.Cells.Item(48, 2).Value = 93.52858140436702
$ws.Cells.Item(48, 3).Value = 192.8312141696829
$ws.Cells.Item(49, 2).Value = 96.01149708713098
$ws.Cells.Item(49, 3).Value = 197.0829639984659

# Append new row 50 with data for x=48, replicating formatting of column A label cells
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = 99.44543991458504
$ws.Cells.Item(50, 3).Value = 200.9484312649789

# Copy the style of A49 (bold, bordered, centered) onto the new A50 label cell
$ws.Cells.Item(49, 1).Copy() | Out-Null
$ws.Cells.Item(50, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
